# Refresh the "Price" (D) and "Volume(1h)" (E) columns with the latest
# crypto market snapshot (GitHub Actions scheduled update).
#
# Price values are stored as plain text in the source data (they use a
# locale-style "." as both thousands- and decimal-separator, e.g.
# "65.182.93"), so each Price cell is briefly switched to the Text number
# format before its value is set -- this stops Excel from reinterpreting
# the string as a number -- and then the cell style is reset back to
# "Normal" so the cell's formatting matches the rest of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.182.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.384.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "184.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -9.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "530.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.603"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.381.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.626"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.37%  "
$ws.Range("E12").Value = "  -12.40%  "
$ws.Range("E13").Value = "  -11.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.921.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.69%  "
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.382.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.967.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -9.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -11.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.967"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "373.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.65%  "
$ws.Range("E24").Value = "  -12.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -16.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.26%  "
$ws.Range("E27").Value = "  -5.19%  "
$ws.Range("E28").Value = "  -10.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "667.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -13.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "60.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.104"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.42%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -13.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.379"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.127"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.823.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -12.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -14.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0625"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -20.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0391"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.57%  "
$ws.Range("E47").Value = "  -13.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.124"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "136.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.94%  "
